$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36 (pushes old rows 36..54 down to 37..55)
$ws.Rows("36:36").Insert()

$ws.Range("A36").Value = "z0bug.invoice_Z0_4_2"
$ws.Range("B36").Value = "z0bug.invoice_Z0_4"
$ws.Range("D36").Value = "z0bug.product_product_26"
$ws.Range("E36").Value = "Prodotto Zeta (versione EU)"
$ws.Range("F36").Value = 100
$ws.Range("G36").Value = "z0bug.coa_510100"
$ws.Range("H36").Value = 1.5
$ws.Range("I36").Value = "z0bug.tax_a41v"
